$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Jumlah" header (column J) to "Jumlah_Transaksi"
$ws.Range("J1").Value = "Jumlah_Transaksi"

# Remove the sample data rows (rows 2-4), keeping only the header row
$ws.Range("A2:J4").EntireRow.Delete()
